# Updated daily measurement data
# Rows 25-27, column G ("par") were changed from the placeholder text "NA"
# to actual measured numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G25").Value = 580
$ws.Range("G26").Value = 582
$ws.Range("G27").Value = 590

# Reflect the active cell/selection left behind after the edit
$null = $ws.Range("G28").Select()
